# "added my first tests and test cases and merged the file with calvins test cases"
#
# This reproduces the content-level changes from the commit:
#  - Testprotokoll (sheet1): the test-case numbering in column A for rows
#    13-19 is cleared (only test cases #1 and #2 remain filled in), and the
#    leftover "highlight" direct formatting in column D for rows 11-17 is
#    removed so it matches the plain look already used further down the
#    table.
#  - The active sheet/selection moves to the Testbericht sheet (Calvin's
#    test report), with the cursor left on cell D13 in Testprotokoll and
#    I8 in Testbericht.
#  - Dependent formulas (e.g. the COUNTA summary on Testbericht) are left
#    as formulas so they recalculate naturally from the new data.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Testprotokoll")
$ws2 = $wb.Worksheets.Item("Testbericht")

# --- Testprotokoll: drop test case numbers 3-9 (rows 13-19), keeping only
#     test cases 1 and 2 filled in; reformat those cells to the plain
#     (no alternating border) look already used by the rows below them.
$ws1.Range("A20").Copy()
$ws1.Range("A13:A19").PasteSpecial(-4122)
$ws1.Range("A13:A19").ClearContents()

# --- Testprotokoll: remove the leftover green "highlight" formatting in
#     column D for rows 11-17 so the cells look like the rest of the table.
$ws1.Range("B11").Copy()
$ws1.Range("D11").PasteSpecial(-4122)
$ws1.Range("D18").Copy()
$ws1.Range("D12:D17").PasteSpecial(-4122)

$excel.CutCopyMode = 0

# --- Selections / active sheet: work finished on the Testbericht sheet.
[void]$ws1.Activate()
[void]$ws1.Range("D13").Select()
[void]$ws2.Activate()
[void]$ws2.Range("I8").Select()
